# Add a "Swiss" test-data worksheet (Switzerland market), cloned from the
# existing "Czech" sheet, matching the layout/styling used by the other
# country sheets in this workbook.

$wb = $excel.ActiveWorkbook

$czech = $wb.Worksheets.Item("Czech")

# Select the source sheet and its full range first so the sheet we copy
# from ends up with a "whole sheet" selection afterwards (same pattern
# used by the other inactive country sheets already in the workbook).
$czech.Select()
$czech.Cells.Select()

# Duplicate "Czech" and place the copy after the last sheet in the workbook.
$czech.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the market name and product code for the new Switzerland sheet.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2645"

# Leave the same cell selected as on the other active/tabbed sheets.
$swiss.Range("A12").Select()
